# Upload new version with timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: LASILACTONE 50/20MG 30 TAB. ---
# H8 (balance) "3:2" -> "3:0" (text, not numeric-looking, plain assign keeps style)
$ws.Range("H8").Value = "3:0"

# P8 (sell price) "41.5800" -> "126.0000" - this looks like a number, so Excel's
# smart entry would coerce it to a real number and drop the t="s" text typing.
# Temporarily switch to a text format, write it, then restore the original
# number format so the stored value stays textual while the cell format is
# unaffected.
$p8 = $ws.Range("P8")
$p8fmt = $p8.NumberFormat()
$p8.NumberFormat = "@"
$p8.Value = "126.0000"
$p8.NumberFormat = $p8fmt

# Q8 (transactions) "0:1" -> "1:0"
$ws.Range("Q8").Value = "1:0"

# --- Row 9: NEXICURE 40 MG 20 F.C. TABLETS ---
# H9 (balance) "0:1" -> "0:0"
$ws.Range("H9").Value = "0:0"

# P9 (sell price) "76.0000" -> "152.0000" (same numeric-looking text issue as P8)
$p9 = $ws.Range("P9")
$p9fmt = $p9.NumberFormat()
$p9.NumberFormat = "@"
$p9.Value = "152.0000"
$p9.NumberFormat = $p9fmt

# Q9 (transactions) "0:1" -> "1:0"
$ws.Range("Q9").Value = "1:0"

# --- Row 10: totals ---
# P10 numeric total 170.58000000000001 -> 331
$ws.Range("P10").Value = 331

# --- Row 11: footer timestamp ---
# A11 "Tuesday, 19 August, 2025 9:57 AM" -> "Tuesday, 19 August, 2025 9:58 AM"
$ws.Range("A11").Value = "Tuesday, 19 August, 2025 9:58 AM"
